# 9.5.1 R&D expenditure table — update to newer edition:
#  - title row loses the trailing "*" (footnote marker moved to a new
#    "preliminary data" column instead)
#  - year columns shift from 2008-2018 to 2017-2020 plus a "2021*" column
#  - data values updated
#  - footnote text replaced with a "preliminary data" note
#  - unused trailing columns (which used to hold 2008-2018) are removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-unused columns for each row (column extent differs by
#     row, matching how the sheet was actually trimmed) ---
$ws.Range("F1:N2").Clear()
$ws.Range("I3:N5").Clear()
$ws.Range("F6:N8").Clear()

# --- Row 1: title (no trailing "*" anymore) ---
$ws.Range("A1").Value = "9.5.1 ИДП га болгон тажрыйбалык-конструктордук жумуштун жана илимий изилдөөнүн чыгымдарынын үлүшү"
$ws.Range("B1").Value = "9.5.1  Доля расходов на научно-исследовательские и опытно-конструкторские работы в ВВП"
$ws.Range("C1").Value = "9.5.1 Research and development expenditure as a proportion of GDP"
$ws.Rows(1).RowHeight = 43.5

# --- Row 4: year headers 2017-2020 plus a "2021*" (preliminary) column ---
$ws.Range("D4").Value = 2017
$ws.Range("E4").Value = 2018
$ws.Range("F4").Value = 2019
$ws.Range("G4").Value = 2020
$ws.Range("H4").Value = "2021*"
$ws.Range("H4").HorizontalAlignment = -4152

# --- Row 5: data values ---
$ws.Range("D5").Value = 0.11
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 0.09
$ws.Range("G5").Value = 0.09
$ws.Range("H5").Value = 0.08

# --- Row 6: footnote, now about preliminary data (kg/ru/en) ---
$ws.Range("A6").Value = "*алдын алаа маалыматтар"
$ws.Range("B6").Value = "*предварительные данные"
$ws.Range("C6").Value = "*preliminary data"
